$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook

# Apply renames/casing fixes to the "Variables" sheet
$wsVars = $wb.Worksheets.Item("Variables")

$wsVars.Range("B2").Value = "id"
$wsVars.Range("B29").Value = "med_stat"
$wsVars.Range("B30").Value = "med_nsaid"
$wsVars.Range("B38").Value = "f1_untdat"
$wsVars.Range("B39").Value = "f2_untdat"

# Apply changes to the "Categories" sheet
$wsCats = $wb.Worksheets.Item("Categories")

$wsCats.Range("C56").Value = "Yes"
$wsCats.Range("C64").Value = "I don't know"

$wsCats.Range("A69").Value = "med_stat"
$wsCats.Range("A70").Value = "med_stat"
$wsCats.Range("A71").Value = "med_nsaid"
$wsCats.Range("A72").Value = "med_nsaid"
